$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D; existing D:K data shifts to E:L.
$ws.Columns("D").Insert()

# Copy number formats from the (shifted) E column into the new D column
# so the new column inherits the same per-row formatting (date row vs data rows).
$ws.Range("E5:E102").Copy() | Out-Null
$ws.Range("D5:D102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$row7 = New-Object 'object[,]' 1,9
$row7[0,0] = 43465
$row7[0,1] = 43100
$row7[0,2] = 42735
$row7[0,3] = 42369
$row7[0,4] = 42004
$row7[0,5] = 41639
$row7[0,6] = 41274
$row7[0,7] = 40908
$row7[0,8] = $null
$ws.Range("D7:L7").Value = $row7

$row8 = New-Object 'object[,]' 1,9
$row8[0,0] = 14800
$row8[0,1] = 8000
$row8[0,2] = 8500
$row8[0,3] = 10900
$row8[0,4] = 11400
$row8[0,5] = 44200
$row8[0,6] = 63000
$row8[0,7] = 32700
$row8[0,8] = $null
$ws.Range("D8:L8").Value = $row8

$row9 = New-Object 'object[,]' 1,9
$row9[0,0] = 10500
$row9[0,1] = 5000
$row9[0,2] = 6300
$row9[0,3] = 8500
$row9[0,4] = 7500
$row9[0,5] = 43200
$row9[0,6] = 49700
$row9[0,7] = 24900
$row9[0,8] = $null
$ws.Range("D9:L9").Value = $row9

$row10 = New-Object 'object[,]' 1,9
$row10[0,0] = 4200
$row10[0,1] = 3000
$row10[0,2] = 2100
$row10[0,3] = 2400
$row10[0,4] = 3900
$row10[0,5] = 1000
$row10[0,6] = 13300
$row10[0,7] = 7800
$row10[0,8] = $null
$ws.Range("D10:L10").Value = $row10

$row11 = New-Object 'object[,]' 1,9
$row11[0,0] = $null
$row11[0,1] = $null
$row11[0,2] = $null
$row11[0,3] = $null
$row11[0,4] = $null
$row11[0,5] = $null
$row11[0,6] = $null
$row11[0,7] = $null
$row11[0,8] = $null
$ws.Range("D11:L11").Value = $row11

$row12 = New-Object 'object[,]' 1,9
$row12[0,0] = "NA"
$row12[0,1] = "NA"
$row12[0,2] = "NA"
$row12[0,3] = "NA"
$row12[0,4] = 400
$row12[0,5] = 500
$row12[0,6] = 500
$row12[0,7] = 400
$row12[0,8] = $null
$ws.Range("D12:L12").Value = $row12

$row13 = New-Object 'object[,]' 1,9
$row13[0,0] = 0
$row13[0,1] = 0
$row13[0,2] = 0
$row13[0,3] = 0
$row13[0,4] = 0
$row13[0,5] = 0
$row13[0,6] = 0
$row13[0,7] = 0
$row13[0,8] = $null
$ws.Range("D13:L13").Value = $row13

$row14 = New-Object 'object[,]' 1,9
$row14[0,0] = 100
$row14[0,1] = 1000
$row14[0,2] = "NA"
$row14[0,3] = 100
$row14[0,4] = 22300
$row14[0,5] = 42800
$row14[0,6] = "NA"
$row14[0,7] = "NA"
$row14[0,8] = $null
$ws.Range("D14:L14").Value = $row14

$row15 = New-Object 'object[,]' 1,9
$row15[0,0] = 1200
$row15[0,1] = 1300
$row15[0,2] = 1400
$row15[0,3] = 1500
$row15[0,4] = 2000
$row15[0,5] = 700
$row15[0,6] = 300
$row15[0,7] = 100
$row15[0,8] = $null
$ws.Range("D15:L15").Value = $row15

$row16 = New-Object 'object[,]' 1,9
$row16[0,0] = $null
$row16[0,1] = $null
$row16[0,2] = $null
$row16[0,3] = $null
$row16[0,4] = $null
$row16[0,5] = $null
$row16[0,6] = $null
$row16[0,7] = $null
$row16[0,8] = $null
$ws.Range("D16:L16").Value = $row16

$row17 = New-Object 'object[,]' 1,9
$row17[0,0] = 21300
$row17[0,1] = 14600
$row17[0,2] = 13400
$row17[0,3] = 16500
$row17[0,4] = 37900
$row17[0,5] = 90500
$row17[0,6] = 62300
$row17[0,7] = 29300
$row17[0,8] = $null
$ws.Range("D17:L17").Value = $row17

$row18 = New-Object 'object[,]' 1,9
$row18[0,0] = -6500
$row18[0,1] = -6600
$row18[0,2] = -4900
$row18[0,3] = -5600
$row18[0,4] = -26500
$row18[0,5] = -46300
$row18[0,6] = 700
$row18[0,7] = 3400
$row18[0,8] = $null
$ws.Range("D18:L18").Value = $row18

$row19 = New-Object 'object[,]' 1,9
$row19[0,0] = $null
$row19[0,1] = $null
$row19[0,2] = $null
$row19[0,3] = $null
$row19[0,4] = $null
$row19[0,5] = $null
$row19[0,6] = $null
$row19[0,7] = $null
$row19[0,8] = $null
$ws.Range("D19:L19").Value = $row19

$row20 = New-Object 'object[,]' 1,9
$row20[0,0] = 1300
$row20[0,1] = 1900
$row20[0,2] = 2000
$row20[0,3] = 700
$row20[0,4] = 0
$row20[0,5] = 100
$row20[0,6] = 5200
$row20[0,7] = 100
$row20[0,8] = $null
$ws.Range("D20:L20").Value = $row20

$row21 = New-Object 'object[,]' 1,9
$row21[0,0] = -4000
$row21[0,1] = -3400
$row21[0,2] = -1400
$row21[0,3] = -3200
$row21[0,4] = -24300
$row21[0,5] = -32300
$row21[0,6] = 24000
$row21[0,7] = 10700
$row21[0,8] = $null
$ws.Range("D21:L21").Value = $row21

$row22 = New-Object 'object[,]' 1,9
$row22[0,0] = 0
$row22[0,1] = 0
$row22[0,2] = 0
$row22[0,3] = 0
$row22[0,4] = 0
$row22[0,5] = 0
$row22[0,6] = 0
$row22[0,7] = 0
$row22[0,8] = $null
$ws.Range("D22:L22").Value = $row22

$row23 = New-Object 'object[,]' 1,9
$row23[0,0] = -5200
$row23[0,1] = -4700
$row23[0,2] = -3000
$row23[0,3] = -4900
$row23[0,4] = -26500
$row23[0,5] = -46200
$row23[0,6] = 5900
$row23[0,7] = 3500
$row23[0,8] = $null
$ws.Range("D23:L23").Value = $row23

$row24 = New-Object 'object[,]' 1,9
$row24[0,0] = 1100
$row24[0,1] = 1100
$row24[0,2] = 300
$row24[0,3] = -800
$row24[0,4] = -2100
$row24[0,5] = -8700
$row24[0,6] = 3400
$row24[0,7] = 2100
$row24[0,8] = $null
$ws.Range("D24:L24").Value = $row24

$row25 = New-Object 'object[,]' 1,9
$row25[0,0] = 0
$row25[0,1] = 0
$row25[0,2] = 0
$row25[0,3] = 0
$row25[0,4] = 0
$row25[0,5] = 0
$row25[0,6] = 0
$row25[0,7] = 0
$row25[0,8] = $null
$ws.Range("D25:L25").Value = $row25

$row26 = New-Object 'object[,]' 1,9
$row26[0,0] = -6300
$row26[0,1] = -5800
$row26[0,2] = -3200
$row26[0,3] = -4100
$row26[0,4] = -24400
$row26[0,5] = -37500
$row26[0,6] = 2500
$row26[0,7] = 1300
$row26[0,8] = $null
$ws.Range("D26:L26").Value = $row26

$row27 = New-Object 'object[,]' 1,9
$row27[0,0] = -6300
$row27[0,1] = -5800
$row27[0,2] = -3200
$row27[0,3] = -4100
$row27[0,4] = -24400
$row27[0,5] = -37500
$row27[0,6] = 2500
$row27[0,7] = 1300
$row27[0,8] = $null
$ws.Range("D27:L27").Value = $row27

$row28 = New-Object 'object[,]' 1,9
$row28[0,0] = 0
$row28[0,1] = 0
$row28[0,2] = 0
$row28[0,3] = 0
$row28[0,4] = 0
$row28[0,5] = 0
$row28[0,6] = 0
$row28[0,7] = 0
$row28[0,8] = $null
$ws.Range("D28:L28").Value = $row28

$row29 = New-Object 'object[,]' 1,9
$row29[0,0] = 0
$row29[0,1] = 0
$row29[0,2] = 0
$row29[0,3] = 0
$row29[0,4] = 0
$row29[0,5] = 0
$row29[0,6] = 0
$row29[0,7] = 0
$row29[0,8] = $null
$ws.Range("D29:L29").Value = $row29

$row30 = New-Object 'object[,]' 1,9
$row30[0,0] = 0
$row30[0,1] = 0
$row30[0,2] = 0
$row30[0,3] = 0
$row30[0,4] = 0
$row30[0,5] = 0
$row30[0,6] = 0
$row30[0,7] = 0
$row30[0,8] = $null
$ws.Range("D30:L30").Value = $row30

$row31 = New-Object 'object[,]' 1,9
$row31[0,0] = 0
$row31[0,1] = 0
$row31[0,2] = 0
$row31[0,3] = 0
$row31[0,4] = 0
$row31[0,5] = 0
$row31[0,6] = 0
$row31[0,7] = 0
$row31[0,8] = $null
$ws.Range("D31:L31").Value = $row31

$row32 = New-Object 'object[,]' 1,9
$row32[0,0] = -1300
$row32[0,1] = -1900
$row32[0,2] = -2000
$row32[0,3] = -700
$row32[0,4] = 0
$row32[0,5] = -100
$row32[0,6] = -5200
$row32[0,7] = -100
$row32[0,8] = $null
$ws.Range("D32:L32").Value = $row32

$row33 = New-Object 'object[,]' 1,9
$row33[0,0] = -6300
$row33[0,1] = -5800
$row33[0,2] = -3200
$row33[0,3] = -4100
$row33[0,4] = -24400
$row33[0,5] = -37500
$row33[0,6] = 2500
$row33[0,7] = 1300
$row33[0,8] = $null
$ws.Range("D33:L33").Value = $row33

$row34 = New-Object 'object[,]' 1,9
$row34[0,0] = 0
$row34[0,1] = 0
$row34[0,2] = 0
$row34[0,3] = 0
$row34[0,4] = 0
$row34[0,5] = 0
$row34[0,6] = 0
$row34[0,7] = 0
$row34[0,8] = $null
$ws.Range("D34:L34").Value = $row34

$row35 = New-Object 'object[,]' 1,9
$row35[0,0] = -6300
$row35[0,1] = -5800
$row35[0,2] = -3200
$row35[0,3] = -4100
$row35[0,4] = -24400
$row35[0,5] = -37500
$row35[0,6] = 2500
$row35[0,7] = 1300
$row35[0,8] = $null
$ws.Range("D35:L35").Value = $row35

$row38 = New-Object 'object[,]' 1,9
$row38[0,0] = 43465
$row38[0,1] = 43100
$row38[0,2] = 42735
$row38[0,3] = 42369
$row38[0,4] = 42004
$row38[0,5] = 41639
$row38[0,6] = 41274
$row38[0,7] = 40908
$row38[0,8] = $null
$ws.Range("D38:L38").Value = $row38

$row39 = New-Object 'object[,]' 1,9
$row39[0,0] = $null
$row39[0,1] = $null
$row39[0,2] = $null
$row39[0,3] = $null
$row39[0,4] = $null
$row39[0,5] = $null
$row39[0,6] = $null
$row39[0,7] = $null
$row39[0,8] = $null
$ws.Range("D39:L39").Value = $row39

$row40 = New-Object 'object[,]' 1,9
$row40[0,0] = $null
$row40[0,1] = $null
$row40[0,2] = $null
$row40[0,3] = $null
$row40[0,4] = $null
$row40[0,5] = $null
$row40[0,6] = $null
$row40[0,7] = $null
$row40[0,8] = $null
$ws.Range("D40:L40").Value = $row40

$row41 = New-Object 'object[,]' 1,9
$row41[0,0] = 2700
$row41[0,1] = 17800
$row41[0,2] = 4800
$row41[0,3] = 4000
$row41[0,4] = 1900
$row41[0,5] = 3600
$row41[0,6] = 1700
$row41[0,7] = 1400
$row41[0,8] = $null
$ws.Range("D41:L41").Value = $row41

$row42 = New-Object 'object[,]' 1,9
$row42[0,0] = 3900
$row42[0,1] = 9400
$row42[0,2] = 11600
$row42[0,3] = 2100
$row42[0,4] = 4500
$row42[0,5] = 2800
$row42[0,6] = 15500
$row42[0,7] = 30700
$row42[0,8] = $null
$ws.Range("D42:L42").Value = $row42

$row43 = New-Object 'object[,]' 1,9
$row43[0,0] = 5100
$row43[0,1] = 4600
$row43[0,2] = 2200
$row43[0,3] = 1900
$row43[0,4] = 2900
$row43[0,5] = 3700
$row43[0,6] = 7300
$row43[0,7] = 8500
$row43[0,8] = $null
$ws.Range("D43:L43").Value = $row43

$row44 = New-Object 'object[,]' 1,9
$row44[0,0] = 600
$row44[0,1] = 500
$row44[0,2] = 100
$row44[0,3] = 100
$row44[0,4] = 100
$row44[0,5] = 3900
$row44[0,6] = 6500
$row44[0,7] = 6600
$row44[0,8] = $null
$ws.Range("D44:L44").Value = $row44

$row45 = New-Object 'object[,]' 1,9
$row45[0,0] = 700
$row45[0,1] = 800
$row45[0,2] = 300
$row45[0,3] = 3300
$row45[0,4] = 2000
$row45[0,5] = 300
$row45[0,6] = 400
$row45[0,7] = 200
$row45[0,8] = $null
$ws.Range("D45:L45").Value = $row45

$row46 = New-Object 'object[,]' 1,9
$row46[0,0] = 13000
$row46[0,1] = 16700
$row46[0,2] = 19000
$row46[0,3] = 11300
$row46[0,4] = 10600
$row46[0,5] = 14300
$row46[0,6] = 31400
$row46[0,7] = 47300
$row46[0,8] = $null
$ws.Range("D46:L46").Value = $row46

$row47 = New-Object 'object[,]' 1,9
$row47[0,0] = 300
$row47[0,1] = 800
$row47[0,2] = "NA"
$row47[0,3] = 300
$row47[0,4] = 400
$row47[0,5] = 400
$row47[0,6] = 2000
$row47[0,7] = 0
$row47[0,8] = $null
$ws.Range("D47:L47").Value = $row47

$row48 = New-Object 'object[,]' 1,9
$row48[0,0] = 72500
$row48[0,1] = 108200
$row48[0,2] = 59400
$row48[0,3] = 57000
$row48[0,4] = 56300
$row48[0,5] = 75700
$row48[0,6] = 117600
$row48[0,7] = 110400
$row48[0,8] = $null
$ws.Range("D48:L48").Value = $row48

$row49 = New-Object 'object[,]' 1,9
$row49[0,0] = 500
$row49[0,1] = 100
$row49[0,2] = 100
$row49[0,3] = 200
$row49[0,4] = 300
$row49[0,5] = 200
$row49[0,6] = 300
$row49[0,7] = 500
$row49[0,8] = $null
$ws.Range("D49:L49").Value = $row49

$row50 = New-Object 'object[,]' 1,9
$row50[0,0] = 0
$row50[0,1] = 0
$row50[0,2] = 0
$row50[0,3] = 0
$row50[0,4] = 0
$row50[0,5] = 0
$row50[0,6] = 0
$row50[0,7] = 0
$row50[0,8] = $null
$ws.Range("D50:L50").Value = $row50

$row51 = New-Object 'object[,]' 1,9
$row51[0,0] = 0
$row51[0,1] = 0
$row51[0,2] = 0
$row51[0,3] = 0
$row51[0,4] = 0
$row51[0,5] = 0
$row51[0,6] = 0
$row51[0,7] = 0
$row51[0,8] = $null
$ws.Range("D51:L51").Value = $row51

$row52 = New-Object 'object[,]' 1,9
$row52[0,0] = 12700
$row52[0,1] = 13300
$row52[0,2] = 9000
$row52[0,3] = 7500
$row52[0,4] = 10700
$row52[0,5] = 7000
$row52[0,6] = 6600
$row52[0,7] = 3700
$row52[0,8] = $null
$ws.Range("D52:L52").Value = $row52

$row53 = New-Object 'object[,]' 1,9
$row53[0,0] = 0
$row53[0,1] = 0
$row53[0,2] = 0
$row53[0,3] = 0
$row53[0,4] = 0
$row53[0,5] = 0
$row53[0,6] = 0
$row53[0,7] = 0
$row53[0,8] = $null
$ws.Range("D53:L53").Value = $row53

$row54 = New-Object 'object[,]' 1,9
$row54[0,0] = 99000
$row54[0,1] = 91000
$row54[0,2] = 87500
$row54[0,3] = 76300
$row54[0,4] = 78300
$row54[0,5] = 97700
$row54[0,6] = 158000
$row54[0,7] = 161800
$row54[0,8] = $null
$ws.Range("D54:L54").Value = $row54

$row55 = New-Object 'object[,]' 1,9
$row55[0,0] = $null
$row55[0,1] = $null
$row55[0,2] = $null
$row55[0,3] = $null
$row55[0,4] = $null
$row55[0,5] = $null
$row55[0,6] = $null
$row55[0,7] = $null
$row55[0,8] = $null
$ws.Range("D55:L55").Value = $row55

$row56 = New-Object 'object[,]' 1,9
$row56[0,0] = $null
$row56[0,1] = $null
$row56[0,2] = $null
$row56[0,3] = $null
$row56[0,4] = $null
$row56[0,5] = $null
$row56[0,6] = $null
$row56[0,7] = $null
$row56[0,8] = $null
$ws.Range("D56:L56").Value = $row56

$row57 = New-Object 'object[,]' 1,9
$row57[0,0] = 2700
$row57[0,1] = 3800
$row57[0,2] = 600
$row57[0,3] = 700
$row57[0,4] = 1000
$row57[0,5] = 800
$row57[0,6] = 7400
$row57[0,7] = 5900
$row57[0,8] = $null
$ws.Range("D57:L57").Value = $row57

$row58 = New-Object 'object[,]' 1,9
$row58[0,0] = 0
$row58[0,1] = 0
$row58[0,2] = 0
$row58[0,3] = 0
$row58[0,4] = 0
$row58[0,5] = 0
$row58[0,6] = 0
$row58[0,7] = 0
$row58[0,8] = $null
$ws.Range("D58:L58").Value = $row58

$row59 = New-Object 'object[,]' 1,9
$row59[0,0] = 3300
$row59[0,1] = 2000
$row59[0,2] = 1000
$row59[0,3] = 1600
$row59[0,4] = 1900
$row59[0,5] = 2100
$row59[0,6] = 4800
$row59[0,7] = 4600
$row59[0,8] = $null
$ws.Range("D59:L59").Value = $row59

$row60 = New-Object 'object[,]' 1,9
$row60[0,0] = 6000
$row60[0,1] = 3100
$row60[0,2] = 1600
$row60[0,3] = 2200
$row60[0,4] = 2800
$row60[0,5] = 2900
$row60[0,6] = 12200
$row60[0,7] = 10500
$row60[0,8] = $null
$ws.Range("D60:L60").Value = $row60

$row61 = New-Object 'object[,]' 1,9
$row61[0,0] = 0
$row61[0,1] = 0
$row61[0,2] = 0
$row61[0,3] = 0
$row61[0,4] = 0
$row61[0,5] = 0
$row61[0,6] = 0
$row61[0,7] = 0
$row61[0,8] = $null
$ws.Range("D61:L61").Value = $row61

$row62 = New-Object 'object[,]' 1,9
$row62[0,0] = 6200
$row62[0,1] = 12900
$row62[0,2] = 18500
$row62[0,3] = 18200
$row62[0,4] = 18100
$row62[0,5] = 19400
$row62[0,6] = 36700
$row62[0,7] = 44600
$row62[0,8] = $null
$ws.Range("D62:L62").Value = $row62

$row63 = New-Object 'object[,]' 1,9
$row63[0,0] = 0
$row63[0,1] = 0
$row63[0,2] = 0
$row63[0,3] = 0
$row63[0,4] = 0
$row63[0,5] = 0
$row63[0,6] = 0
$row63[0,7] = 0
$row63[0,8] = $null
$ws.Range("D63:L63").Value = $row63

$row64 = New-Object 'object[,]' 1,9
$row64[0,0] = 0
$row64[0,1] = 0
$row64[0,2] = 0
$row64[0,3] = 0
$row64[0,4] = 0
$row64[0,5] = 0
$row64[0,6] = 0
$row64[0,7] = 0
$row64[0,8] = $null
$ws.Range("D64:L64").Value = $row64

$row65 = New-Object 'object[,]' 1,9
$row65[0,0] = 0
$row65[0,1] = 0
$row65[0,2] = 0
$row65[0,3] = 0
$row65[0,4] = 0
$row65[0,5] = 0
$row65[0,6] = 0
$row65[0,7] = 0
$row65[0,8] = $null
$ws.Range("D65:L65").Value = $row65

$row66 = New-Object 'object[,]' 1,9
$row66[0,0] = 12200
$row66[0,1] = 7300
$row66[0,2] = 20100
$row66[0,3] = 20500
$row66[0,4] = 21000
$row66[0,5] = 22400
$row66[0,6] = 48900
$row66[0,7] = 55000
$row66[0,8] = $null
$ws.Range("D66:L66").Value = $row66

$row67 = New-Object 'object[,]' 1,9
$row67[0,0] = $null
$row67[0,1] = $null
$row67[0,2] = $null
$row67[0,3] = $null
$row67[0,4] = $null
$row67[0,5] = $null
$row67[0,6] = $null
$row67[0,7] = $null
$row67[0,8] = $null
$ws.Range("D67:L67").Value = $row67

$row68 = New-Object 'object[,]' 1,9
$row68[0,0] = 0
$row68[0,1] = 0
$row68[0,2] = 0
$row68[0,3] = 0
$row68[0,4] = 0
$row68[0,5] = 0
$row68[0,6] = 0
$row68[0,7] = 0
$row68[0,8] = $null
$ws.Range("D68:L68").Value = $row68

$row69 = New-Object 'object[,]' 1,9
$row69[0,0] = 0
$row69[0,1] = 0
$row69[0,2] = 0
$row69[0,3] = 0
$row69[0,4] = 0
$row69[0,5] = 0
$row69[0,6] = 0
$row69[0,7] = 0
$row69[0,8] = $null
$ws.Range("D69:L69").Value = $row69

$row70 = New-Object 'object[,]' 1,9
$row70[0,0] = 0
$row70[0,1] = 0
$row70[0,2] = 0
$row70[0,3] = 0
$row70[0,4] = 0
$row70[0,5] = 0
$row70[0,6] = 0
$row70[0,7] = 0
$row70[0,8] = $null
$ws.Range("D70:L70").Value = $row70

$row71 = New-Object 'object[,]' 1,9
$row71[0,0] = 0
$row71[0,1] = 0
$row71[0,2] = 0
$row71[0,3] = 0
$row71[0,4] = 0
$row71[0,5] = 0
$row71[0,6] = 0
$row71[0,7] = 0
$row71[0,8] = $null
$ws.Range("D71:L71").Value = $row71

$row72 = New-Object 'object[,]' 1,9
$row72[0,0] = -84400
$row72[0,1] = -77800
$row72[0,2] = -80900
$row72[0,3] = -78100
$row72[0,4] = -73200
$row72[0,5] = -47900
$row72[0,6] = -10300
$row72[0,7] = -15200
$row72[0,8] = $null
$ws.Range("D72:L72").Value = $row72

$row73 = New-Object 'object[,]' 1,9
$row73[0,0] = 0
$row73[0,1] = 0
$row73[0,2] = 0
$row73[0,3] = 0
$row73[0,4] = 0
$row73[0,5] = 0
$row73[0,6] = 0
$row73[0,7] = 0
$row73[0,8] = $null
$ws.Range("D73:L73").Value = $row73

$row74 = New-Object 'object[,]' 1,9
$row74[0,0] = 0
$row74[0,1] = 0
$row74[0,2] = 0
$row74[0,3] = 0
$row74[0,4] = 0
$row74[0,5] = 0
$row74[0,6] = 0
$row74[0,7] = 0
$row74[0,8] = $null
$ws.Range("D74:L74").Value = $row74

$row75 = New-Object 'object[,]' 1,9
$row75[0,0] = 0
$row75[0,1] = 0
$row75[0,2] = 0
$row75[0,3] = 0
$row75[0,4] = 0
$row75[0,5] = 0
$row75[0,6] = 0
$row75[0,7] = 0
$row75[0,8] = $null
$ws.Range("D75:L75").Value = $row75

$row76 = New-Object 'object[,]' 1,9
$row76[0,0] = 86800
$row76[0,1] = 83700
$row76[0,2] = 67500
$row76[0,3] = 55800
$row76[0,4] = 57300
$row76[0,5] = 75300
$row76[0,6] = 109100
$row76[0,7] = 106800
$row76[0,8] = $null
$ws.Range("D76:L76").Value = $row76

$row77 = New-Object 'object[,]' 1,9
$row77[0,0] = 0
$row77[0,1] = 0
$row77[0,2] = 0
$row77[0,3] = 0
$row77[0,4] = 0
$row77[0,5] = 0
$row77[0,6] = 0
$row77[0,7] = 0
$row77[0,8] = $null
$ws.Range("D77:L77").Value = $row77

$row80 = New-Object 'object[,]' 1,9
$row80[0,0] = 43465
$row80[0,1] = 43100
$row80[0,2] = 42735
$row80[0,3] = 42369
$row80[0,4] = 42004
$row80[0,5] = 41639
$row80[0,6] = 41274
$row80[0,7] = 40908
$row80[0,8] = $null
$ws.Range("D80:L80").Value = $row80

$row81 = New-Object 'object[,]' 1,9
$row81[0,0] = -6300
$row81[0,1] = -5800
$row81[0,2] = -3200
$row81[0,3] = -4100
$row81[0,4] = -24400
$row81[0,5] = -37500
$row81[0,6] = 2500
$row81[0,7] = 1300
$row81[0,8] = $null
$ws.Range("D81:L81").Value = $row81

$row82 = New-Object 'object[,]' 1,9
$row82[0,0] = $null
$row82[0,1] = $null
$row82[0,2] = $null
$row82[0,3] = $null
$row82[0,4] = $null
$row82[0,5] = $null
$row82[0,6] = $null
$row82[0,7] = $null
$row82[0,8] = $null
$ws.Range("D82:L82").Value = $row82

$row83 = New-Object 'object[,]' 1,9
$row83[0,0] = 1200
$row83[0,1] = 1300
$row83[0,2] = 1500
$row83[0,3] = 1700
$row83[0,4] = 2200
$row83[0,5] = 13900
$row83[0,6] = 18000
$row83[0,7] = 7200
$row83[0,8] = $null
$ws.Range("D83:L83").Value = $row83

$row84 = New-Object 'object[,]' 1,9
$row84[0,0] = 0
$row84[0,1] = 0
$row84[0,2] = 0
$row84[0,3] = 0
$row84[0,4] = 0
$row84[0,5] = 0
$row84[0,6] = 0
$row84[0,7] = 0
$row84[0,8] = $null
$ws.Range("D84:L84").Value = $row84

$row85 = New-Object 'object[,]' 1,9
$row85[0,0] = 0
$row85[0,1] = 0
$row85[0,2] = 0
$row85[0,3] = 0
$row85[0,4] = 0
$row85[0,5] = 0
$row85[0,6] = 0
$row85[0,7] = 0
$row85[0,8] = $null
$ws.Range("D85:L85").Value = $row85

$row86 = New-Object 'object[,]' 1,9
$row86[0,0] = 0
$row86[0,1] = 0
$row86[0,2] = 0
$row86[0,3] = 0
$row86[0,4] = 0
$row86[0,5] = 0
$row86[0,6] = 0
$row86[0,7] = 0
$row86[0,8] = $null
$ws.Range("D86:L86").Value = $row86

$row87 = New-Object 'object[,]' 1,9
$row87[0,0] = 0
$row87[0,1] = 0
$row87[0,2] = 0
$row87[0,3] = 0
$row87[0,4] = 0
$row87[0,5] = 0
$row87[0,6] = 0
$row87[0,7] = 0
$row87[0,8] = $null
$ws.Range("D87:L87").Value = $row87

$row88 = New-Object 'object[,]' 1,9
$row88[0,0] = 0
$row88[0,1] = 0
$row88[0,2] = 0
$row88[0,3] = 0
$row88[0,4] = 0
$row88[0,5] = 0
$row88[0,6] = 0
$row88[0,7] = 0
$row88[0,8] = $null
$ws.Range("D88:L88").Value = $row88

$row89 = New-Object 'object[,]' 1,9
$row89[0,0] = -4100
$row89[0,1] = -3000
$row89[0,2] = -3400
$row89[0,3] = -2000
$row89[0,4] = -500
$row89[0,5] = 2500
$row89[0,6] = 10500
$row89[0,7] = 4300
$row89[0,8] = $null
$ws.Range("D89:L89").Value = $row89

$row90 = New-Object 'object[,]' 1,9
$row90[0,0] = $null
$row90[0,1] = $null
$row90[0,2] = $null
$row90[0,3] = $null
$row90[0,4] = $null
$row90[0,5] = $null
$row90[0,6] = $null
$row90[0,7] = $null
$row90[0,8] = $null
$ws.Range("D90:L90").Value = $row90

$row91 = New-Object 'object[,]' 1,9
$row91[0,0] = -13100
$row91[0,1] = -12100
$row91[0,2] = -4000
$row91[0,3] = -1500
$row91[0,4] = -4800
$row91[0,5] = -16500
$row91[0,6] = -24900
$row91[0,7] = -10500
$row91[0,8] = $null
$ws.Range("D91:L91").Value = $row91

$row92 = New-Object 'object[,]' 1,9
$row92[0,0] = 0
$row92[0,1] = 0
$row92[0,2] = 0
$row92[0,3] = 0
$row92[0,4] = 0
$row92[0,5] = 0
$row92[0,6] = 0
$row92[0,7] = 0
$row92[0,8] = $null
$ws.Range("D92:L92").Value = $row92

$row93 = New-Object 'object[,]' 1,9
$row93[0,0] = 0
$row93[0,1] = 0
$row93[0,2] = 0
$row93[0,3] = 0
$row93[0,4] = 0
$row93[0,5] = 0
$row93[0,6] = 0
$row93[0,7] = 0
$row93[0,8] = $null
$ws.Range("D93:L93").Value = $row93

$row94 = New-Object 'object[,]' 1,9
$row94[0,0] = -10500
$row94[0,1] = -5500
$row94[0,2] = -1300
$row94[0,3] = -1000
$row94[0,4] = -4800
$row94[0,5] = -16800
$row94[0,6] = -24800
$row94[0,7] = -11200
$row94[0,8] = $null
$ws.Range("D94:L94").Value = $row94

$row95 = New-Object 'object[,]' 1,9
$row95[0,0] = $null
$row95[0,1] = $null
$row95[0,2] = $null
$row95[0,3] = $null
$row95[0,4] = $null
$row95[0,5] = $null
$row95[0,6] = $null
$row95[0,7] = $null
$row95[0,8] = $null
$ws.Range("D95:L95").Value = $row95

$row96 = New-Object 'object[,]' 1,9
$row96[0,0] = 0
$row96[0,1] = 0
$row96[0,2] = 0
$row96[0,3] = 0
$row96[0,4] = 0
$row96[0,5] = 0
$row96[0,6] = 0
$row96[0,7] = 0
$row96[0,8] = $null
$ws.Range("D96:L96").Value = $row96

$row97 = New-Object 'object[,]' 1,9
$row97[0,0] = 0
$row97[0,1] = 0
$row97[0,2] = 0
$row97[0,3] = 0
$row97[0,4] = 0
$row97[0,5] = 0
$row97[0,6] = 0
$row97[0,7] = 0
$row97[0,8] = $null
$ws.Range("D97:L97").Value = $row97

$row98 = New-Object 'object[,]' 1,9
$row98[0,0] = 0
$row98[0,1] = 0
$row98[0,2] = 0
$row98[0,3] = 0
$row98[0,4] = 0
$row98[0,5] = 0
$row98[0,6] = 0
$row98[0,7] = 0
$row98[0,8] = $null
$ws.Range("D98:L98").Value = $row98

$row99 = New-Object 'object[,]' 1,9
$row99[0,0] = 0
$row99[0,1] = 0
$row99[0,2] = 0
$row99[0,3] = 0
$row99[0,4] = 0
$row99[0,5] = 0
$row99[0,6] = 0
$row99[0,7] = 0
$row99[0,8] = $null
$ws.Range("D99:L99").Value = $row99

$row100 = New-Object 'object[,]' 1,9
$row100[0,0] = 7700
$row100[0,1] = 6600
$row100[0,2] = 13800
$row100[0,3] = 2700
$row100[0,4] = 5300
$row100[0,5] = 3500
$row100[0,6] = 400
$row100[0,7] = 200
$row100[0,8] = $null
$ws.Range("D100:L100").Value = $row100

$row101 = New-Object 'object[,]' 1,9
$row101[0,0] = 0
$row101[0,1] = 0
$row101[0,2] = 0
$row101[0,3] = 0
$row101[0,4] = 0
$row101[0,5] = 0
$row101[0,6] = 0
$row101[0,7] = 0
$row101[0,8] = $null
$ws.Range("D101:L101").Value = $row101

$row102 = New-Object 'object[,]' 1,9
$row102[0,0] = -6900
$row102[0,1] = -1800
$row102[0,2] = 9100
$row102[0,3] = -400
$row102[0,4] = 0
$row102[0,5] = -10800
$row102[0,6] = -13900
$row102[0,7] = -6700
$row102[0,8] = $null
$ws.Range("D102:L102").Value = $row102
